$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "Förändrad" (changed) date column C for rows 2-6 from 45204 to 45207
foreach ($row in 2..6) {
    $ws.Cells.Item($row, 3).Value = 45207
}
